# Generate Report for Handoff
# Adds two new localization entries (099285d2-... and c6b0a074-...) to the
# existing "ed1873f5-..." / "27e84707-..." rows on all three worksheets
# (Overview, zh-cn, de-de), keeping the existing alphabetical-by-guid order:
#   ed1873f5 (row2) < 099285d2 (new row3) < 27e84707 (row4, was row3) < c6b0a074 (new row5)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"  (columns: File Name | zh-cn | de-de | Latest Handoff Date)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "099285d2-3204-431f-8a21-7bf1d87f9cde.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-21 12:32:40"

$ws.Range("A5").Value = "c6b0a074-04c6-448a-990d-f06b11ecdcb6.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-21 12:32:40"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f1fbdc764a1eb1728d78bddacd5a2eb522d87864/e2e/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md", "", "", "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/099285d2099285d2099285d2099285d2099285d2/e2e/099285d2-3204-431f-8a21-7bf1d87f9cde.md", "", "", "099285d2-3204-431f-8a21-7bf1d87f9cde.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/279f5661bdb22622bb15671b69de9b9388e73433/e2e/27e84707-347a-4dd8-97db-1852d3153c02.md", "", "", "27e84707-347a-4dd8-97db-1852d3153c02.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6b0a074c6b0a074c6b0a074c6b0a074c6b0a074/e2e/c6b0a074-04c6-448a-990d-f06b11ecdcb6.md", "", "", "c6b0a074-04c6-448a-990d-f06b11ecdcb6.md") | Out-Null

# ---------------------------------------------------------------------------
# Shared column headers / values used on the two detail sheets
# ---------------------------------------------------------------------------
function Fill-DetailSheet {
    param($ws, $langCode)

    $ws.Hyperlinks.Delete()

    $ws.Rows.Item(3).Insert()

    # Row-insert copies the formatting pattern (incl. empty F/G cells) from
    # the row above; the new row doesn't use those columns, so drop them.
    $ws.Range("F3").Clear()
    $ws.Range("G3").Clear()

    # --- row 3 : 099285d2 (new) ---
    $ws.Range("A3").Value = "099285d2-3204-431f-8a21-7bf1d87f9cde.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = "099285d2-3204-431f-8a21-7bf1d87f9cde.87f16fdc5007fc241071e2b4ebf5ea50674d3dd4.$langCode.xlf"
    if ($langCode -eq "zh-cn") {
        $ws.Range("E3").Value = "2016-03-21 12:32:37"
    } else {
        $ws.Range("E3").Value = "2016-03-21 12:32:40"
    }
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("J3").Value = "Include"

    # --- row 4 : 27e84707 (shifted down from old row 3, same text values) ---
    $ws.Range("A4").Value = "27e84707-347a-4dd8-97db-1852d3153c02.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = "27e84707-347a-4dd8-97db-1852d3153c02.658a352af91171a2bea34b4812c3d7a6266d09ba.$langCode.xlf"
    if ($langCode -eq "zh-cn") {
        $ws.Range("E4").Value = "2016-03-21 12:30:55"
    } else {
        $ws.Range("E4").Value = "2016-03-21 12:30:59"
    }
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("J4").Value = "Include"

    # --- row 5 : c6b0a074 (new, appended) ---
    $ws.Range("A5").Value = "c6b0a074-04c6-448a-990d-f06b11ecdcb6.md"
    $ws.Range("B5").Value = ".md"
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Range("D5").Value = "c6b0a074-04c6-448a-990d-f06b11ecdcb6.bb5bde3b0188c7d562cef42b2235d302854be512.$langCode.xlf"
    if ($langCode -eq "zh-cn") {
        $ws.Range("E5").Value = "2016-03-21 12:32:37"
    } else {
        $ws.Range("E5").Value = "2016-03-21 12:32:40"
    }
    $ws.Range("H5").Value = "0001-01-01 00:00:00"
    $ws.Range("J5").Value = "Include"

    # --- hyperlinks (A2/D2 keep pointing at the first, untouched row) ---
    $edBase = "https://github.com/OpenLocalizationTest/oltest/blob/f1fbdc764a1eb1728d78bddacd5a2eb522d87864/e2e/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md"
    $ws.Hyperlinks.Add($ws.Range("A2"), $edBase, "", "", "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md") | Out-Null

    if ($langCode -eq "zh-cn") {
        $edXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b97ca0b0406dc0d90487f7d5e82aea599c8e68f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.zh-cn.xlf"
        $edMdTarget  = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ab24cfb7024c90191185b92a45828117cc555d51/e2e/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md"
        $edBackTarget = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9867fb602138c6bc7c773807e52b46533d69e05e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.zh-cn.xlf"
        $newHash099 = "87f16fdc5007fc241071e2b4ebf5ea50674d3dd4"
        $twoSevenTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1b99df7ac53edd523059197298274c2561cb317/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/27e84707-347a-4dd8-97db-1852d3153c02.658a352af91171a2bea34b4812c3d7a6266d09ba.zh-cn.xlf"
    } else {
        $edXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0c2a58e9ed94a5aaf4ee58bd3da4b1cb9b361e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.de-de.xlf"
        $edMdTarget  = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6e10c094d3610c2f4706f2ba8dd7382cd3a708bc/e2e/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md"
        $edBackTarget = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f47ba1e60b26cee0f25bbf88273c13636ac0cd68/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.de-de.xlf"
        $newHash099 = "87f16fdc5007fc241071e2b4ebf5ea50674d3dd4"
        $twoSevenTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3741e64e1eb1508f7a9c969c3789c7e196e81b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/27e84707-347a-4dd8-97db-1852d3153c02.658a352af91171a2bea34b4812c3d7a6266d09ba.de-de.xlf"
    }

    $ws.Hyperlinks.Add($ws.Range("D2"), $edXlfTarget, "", "", "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.$langCode.xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $edBase, "", "", "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $edXlfTarget, "", "", "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.$langCode.xlf") | Out-Null

    $zeroNineNineMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/099285d2099285d2099285d2099285d2099285d2/e2e/099285d2-3204-431f-8a21-7bf1d87f9cde.md"
    $zeroNineNineXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash099/ol-handoff/OpenLocalizationTestOrg/oltest.$langCode/ci/ht/099285d2-3204-431f-8a21-7bf1d87f9cde.87f16fdc5007fc241071e2b4ebf5ea50674d3dd4.$langCode.xlf"
    $ws.Hyperlinks.Add($ws.Range("A3"), $zeroNineNineMdTarget, "", "", "099285d2-3204-431f-8a21-7bf1d87f9cde.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $zeroNineNineXlfTarget, "", "", "099285d2-3204-431f-8a21-7bf1d87f9cde.87f16fdc5007fc241071e2b4ebf5ea50674d3dd4.$langCode.xlf") | Out-Null

    $twoSevenMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/279f5661bdb22622bb15671b69de9b9388e73433/e2e/27e84707-347a-4dd8-97db-1852d3153c02.md"
    $ws.Hyperlinks.Add($ws.Range("A4"), $twoSevenMdTarget, "", "", "27e84707-347a-4dd8-97db-1852d3153c02.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D4"), $twoSevenTarget, "", "", "27e84707-347a-4dd8-97db-1852d3153c02.658a352af91171a2bea34b4812c3d7a6266d09ba.$langCode.xlf") | Out-Null

    $c6MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/c6b0a074c6b0a074c6b0a074c6b0a074c6b0a074/e2e/c6b0a074-04c6-448a-990d-f06b11ecdcb6.md"
    $c6XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb5bde3b0188c7d562cef42b2235d302854be512/ol-handoff/OpenLocalizationTestOrg/oltest.$langCode/ci/ht/c6b0a074-04c6-448a-990d-f06b11ecdcb6.bb5bde3b0188c7d562cef42b2235d302854be512.$langCode.xlf"
    $ws.Hyperlinks.Add($ws.Range("A5"), $c6MdTarget, "", "", "c6b0a074-04c6-448a-990d-f06b11ecdcb6.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D5"), $c6XlfTarget, "", "", "c6b0a074-04c6-448a-990d-f06b11ecdcb6.bb5bde3b0188c7d562cef42b2235d302854be512.$langCode.xlf") | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Fill-DetailSheet $wsZh "zh-cn"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Fill-DetailSheet $wsDe "de-de"

Write-Host "Report generated for handoff."
